# "Updated data from dining services"
#
# Sheet 1 - Compost Feb & April:
#   - The "meal number" column (G) for the Feb data block (rows 8-86) is
#     updated in place: each value advances by one meal in the 1-28 cycle
#     (28 wraps back around to 1). Columns A-F are untouched.
#   - The final row (87, date 2020-05-04 / meal 28) is removed entirely,
#     so the table now ends at row 86.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1 - Compost Feb & April")

# Advance the meal-number cycle (1-28) by one step for every data row
# from row 8 through row 86.
for ($r = 8; $r -le 86; $r++) {
    $current = $ws.Cells.Item($r, 7).Value()
    $next = $current + 1
    if ($next -gt 28) {
        $next = 1
    }
    $ws.Cells.Item($r, 7).Value = $next
}

# Drop the now-removed last day of data (old row 87).
$ws.Rows(87).Delete()
